$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new log entry row (row 73) with text, hours and date.
$ws.Range("B73").Value = "Psaní - formulare, otazky, role, ...; reseni tabulek"
$ws.Range("C73").Value = 3
$ws.Range("D73").Value = 41009

# Make sure the SUM formula in C3 has recalculated before we touch
# formatting (copy/paste can otherwise leave the cached value stale).
$wb.Application.Calculate()

# Copy the row-73 formatting (date style etc.) from the row above it,
# same as dragging the row format down while filling in a new entry.
$ws.Range("B72:D72").Copy()
$ws.Range("B73").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the selection where Excel would land after tabbing through the
# new row (one cell past the last filled-in cell).
$ws.Range("D74").Select()
